$wb = $excel.ActiveWorkbook

# Insert the new "Adjusted2" worksheet right after the "Adjusted" sheet,
# mirroring the other data sheets (measured / calculated / Adjusted).
$adjusted = $wb.Worksheets.Item("Adjusted")
$ws = $wb.Worksheets.Add($null, $adjusted)
$ws.Name = "Adjusted2"

# Match the outline settings used on the other sheets (summary rows below /
# summary columns to the right of the detail).
$ws.Outline.SummaryRow = 1
$ws.Outline.SummaryColumn = 1

# Match the page margins used throughout the rest of the workbook.
$ws.PageSetup.LeftMargin = $excel.InchesToPoints(0.75)
$ws.PageSetup.RightMargin = $excel.InchesToPoints(0.75)
$ws.PageSetup.TopMargin = $excel.InchesToPoints(1)
$ws.PageSetup.BottomMargin = $excel.InchesToPoints(1)
$ws.PageSetup.HeaderMargin = $excel.InchesToPoints(0.5)
$ws.PageSetup.FooterMargin = $excel.InchesToPoints(0.5)

# Header row
$headers = @("theta", "Jxx", "Jyy", "beta", "gamma", "trace")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# Match the bold/centered/bordered header style already used on the other sheets
# by copying the formatting from the "Adjusted" sheet's header row.
$adjusted.Range("A1:F1").Copy()
$ws.Range("A1:F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Data rows (theta, Jxx, Jyy, beta, gamma, trace)
$data = New-Object 'object[,]' 19,6
$data[0,0] = 0
$data[0,1] = 0.0006887052346870647
$data[0,2] = 0.999311301418417
$data[0,3] = 0.001386708065758435
$data[0,4] = 0
$data[0,5] = 0.9986273973759898
$data[1,0] = 5
$data[1,1] = 0.02716108452717102
$data[1,2] = 0.9732370424326104
$data[1,3] = -0.1625859059589901
$data[1,4] = 0
$data[1,5] = 1.000796418908678
$data[2,0] = 10
$data[2,1] = 0.1149122167396198
$data[2,2] = 0.8857955153245778
$data[2,3] = -0.3190434661154541
$data[2,4] = 0
$data[2,5] = 1.001415979067074
$data[3,0] = 15
$data[3,1] = 0.2491450626833772
$data[3,2] = 0.7517097249621988
$data[3,3] = -0.4327641065136215
$data[3,4] = 0
$data[3,5] = 1.001710316635315
$data[4,0] = 20
$data[4,1] = 0.4010155596544819
$data[4,2] = 0.5998933824633242
$data[4,3] = -0.4904758735584562
$data[4,4] = 0
$data[4,5] = 1.001818714494147
$data[5,0] = 25
$data[5,1] = 0.5802814408440432
$data[5,2] = 0.4203885644770028
$data[5,3] = -0.4939065535356619
$data[5,4] = 0
$data[5,5] = 1.001340462982025
$data[6,0] = 30
$data[6,1] = 0.7469675811998441
$data[6,2] = 0.2530323967879994
$data[6,3] = -0.4346361159641195
$data[6,4] = 0
$data[6,5] = 0.9998030677885762
$data[7,0] = 35
$data[7,1] = 0.8749999522573777
$data[7,2] = 0.1250000241203859
$data[7,3] = -0.3288904983564416
$data[7,4] = 0
$data[7,5] = 0.9975878422988075
$data[8,0] = 40
$data[8,1] = 0.9635599038987996
$data[8,2] = 0.03644008401447181
$data[8,3] = -0.178696574538335
$data[8,4] = 0
$data[8,5] = 0.9936404996279151
$data[9,0] = 45
$data[9,1] = 0.9950738935951535
$data[9,2] = 0.004926108353979964
$data[9,3] = -0.007874032793694044
$data[9,4] = 0
$data[9,5] = 0.9903203210430065
$data[10,0] = 50
$data[10,1] = 0.9684844083951484
$data[10,2] = 0.03151558018545603
$data[10,3] = 0.1600566648067495
$data[10,4] = 0
$data[10,5] = 0.9901915529970468
$data[11,0] = 55
$data[11,1] = 0.8861031058426012
$data[11,2] = 0.1138968520728944
$data[11,3] = 0.3116045934683438
$data[11,4] = 0
$data[11,5] = 0.9923460524371626
$data[12,0] = 60
$data[12,1] = 0.7598707014998497
$data[12,2] = 0.2401292277364912
$data[12,3] = 0.4235463032925153
$data[12,4] = 0
$data[12,5] = 0.993848471076708
$data[13,0] = 65
$data[13,1] = 0.5938173243659502
$data[13,2] = 0.4061825851636688
$data[13,3] = 0.4884974709780901
$data[13,4] = 0
$data[13,5] = 0.9948628655113572
$data[14,0] = 70
$data[14,1] = 0.4101123399312268
$data[14,2] = 0.5898875674344477
$data[14,3] = 0.4891151555979193
$data[14,4] = 0
$data[14,5] = 0.9946267444487499
$data[15,0] = 75
$data[15,1] = 0.2627478908029525
$data[15,2] = 0.7372520424212111
$data[15,3] = 0.4373229442429809
$data[15,4] = 0
$data[15,5] = 0.9950797432983463
$data[16,0] = 80
$data[16,1] = 0.1273314383076676
$data[16,2] = 0.872668530423148
$data[16,3] = 0.330703021292808
$data[16,4] = 0
$data[16,5] = 0.9964926357567788
$data[17,0] = 85
$data[17,1] = 0.03865607235227576
$data[17,2] = 0.9613439187653205
$data[17,3] = 0.1889450954510491
$data[17,4] = 0
$data[17,5] = 0.9970769202667797
$data[18,0] = 90
$data[18,1] = 0.001078360891874376
$data[18,2] = 0.9989216454069046
$data[18,3] = 0.0111966217542146
$data[18,4] = 0
$data[18,5] = 0.9980963452020648
$ws.Range("A2:F20").Value = $data

$ws.Range("A1").Select()
